# Auto-generated update of Ultima Profits market-price / leve-profit values
# across the ALC / ARM / BSM / CRP / CUL / GSM / LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3174.875
$ws.Range("I18").Value = 1939.8
$ws.Range("J18").Value = 5233.3335
$ws.Range("K18").Value = 1939.8
$ws.Range("L18").Value = 5233.3335
$ws.Range("M18").Value = -1655.8
$ws.Range("N18").Value = -5801.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 37037420
$ws.Range("I33").Value = 52632084
$ws.Range("J33").Value = 99.75
$ws.Range("K33").Value = 52632084
$ws.Range("L33").Value = 99.75
$ws.Range("M33").Value = -52631855
$ws.Range("N33").Value = -557.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1162.0625
$ws.Range("I41").Value = 1306.2307
$ws.Range("K41").Value = 1306.2307
$ws.Range("M41").Value = -866.2307000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1985.2307
$ws.Range("I116").Value = 2038.75
$ws.Range("J116").Value = 1899.6
$ws.Range("K116").Value = 2038.75
$ws.Range("L116").Value = 1899.6
$ws.Range("M116").Value = 1403.25
$ws.Range("N116").Value = -8783.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9335.690000000001
$ws.Range("I32").Value = 9879.8125
$ws.Range("J32").Value = 5604.5713
$ws.Range("K32").Value = 9879.8125
$ws.Range("L32").Value = 5604.5713
$ws.Range("M32").Value = -9592.8125
$ws.Range("N32").Value = -6178.5713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2121.4285
$ws.Range("I45").Value = 2130
$ws.Range("J45").Value = 2100
$ws.Range("K45").Value = 2130
$ws.Range("L45").Value = 2100
$ws.Range("M45").Value = -1753
$ws.Range("N45").Value = -2854

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1506.3636
$ws.Range("I110").Value = 1413.3334
$ws.Range("J110").Value = 1925
$ws.Range("K110").Value = 1413.3334
$ws.Range("L110").Value = 1925
$ws.Range("M110").Value = 631.6666
$ws.Range("N110").Value = -6015

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6117.879
$ws.Range("I122").Value = 6359.115
$ws.Range("J122").Value = 5221.857
$ws.Range("K122").Value = 19077.345
$ws.Range("L122").Value = 15665.571
$ws.Range("M122").Value = -16627.345
$ws.Range("N122").Value = -20565.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1779.1666
$ws.Range("I107").Value = 1797.4348
$ws.Range("J107").Value = 1719.1428
$ws.Range("K107").Value = 1797.4348
$ws.Range("L107").Value = 1719.1428
$ws.Range("M107").Value = 122.5652
$ws.Range("N107").Value = -5559.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3306.5588
$ws.Range("I134").Value = 1994.591
$ws.Range("K134").Value = 5983.772999999999
$ws.Range("M134").Value = -3448.772999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 885.25
$ws.Range("I16").Value = 762.125
$ws.Range("J16").Value = 1049.4166
$ws.Range("K16").Value = 762.125
$ws.Range("L16").Value = 1049.4166
$ws.Range("M16").Value = -475.125
$ws.Range("N16").Value = -1623.4166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1223.7727
$ws.Range("I99").Value = 1055.5834
$ws.Range("J99").Value = 1425.6
$ws.Range("K99").Value = 1055.5834
$ws.Range("L99").Value = 1425.6
$ws.Range("M99").Value = 442.4166
$ws.Range("N99").Value = -4421.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 885.25
$ws.Range("I113").Value = 762.125
$ws.Range("J113").Value = 1049.4166
$ws.Range("K113").Value = 762.125
$ws.Range("L113").Value = 1049.4166
$ws.Range("M113").Value = 1407.875
$ws.Range("N113").Value = -5389.4166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1505.1034
$ws.Range("I122").Value = 1456.85
$ws.Range("J122").Value = 1612.3334
$ws.Range("K122").Value = 4370.549999999999
$ws.Range("L122").Value = 4837.0002
$ws.Range("M122").Value = -1920.549999999999
$ws.Range("N122").Value = -9737.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1223.7727
$ws.Range("I126").Value = 1055.5834
$ws.Range("J126").Value = 1425.6
$ws.Range("K126").Value = 3166.7502
$ws.Range("L126").Value = 4276.799999999999
$ws.Range("M126").Value = -696.7501999999999
$ws.Range("N126").Value = -9216.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3217.2727
$ws.Range("J80").Value = 3450
$ws.Range("L80").Value = 10350
$ws.Range("N80").Value = -12222

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 3217.2727
$ws.Range("J83").Value = 3450
$ws.Range("L83").Value = 31050
$ws.Range("N83").Value = -40410

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1095.75
$ws.Range("J92").Value = 1095.75
$ws.Range("L92").Value = 3287.25
$ws.Range("N92").Value = -5783.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 344.2857
$ws.Range("I98").Value = 453.33334
$ws.Range("J98").Value = 262.5
$ws.Range("K98").Value = 1360.00002
$ws.Range("L98").Value = 787.5
$ws.Range("M98").Value = 137.9999800000001
$ws.Range("N98").Value = -3783.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1163.1818
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1163.1818
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 3489.5454
$ws.Range("N107").Value = -7329.5454
$ws.Range("M107").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2450.7666
$ws.Range("I109").Value = 879.36365
$ws.Range("J109").Value = 3360.5264
$ws.Range("K109").Value = 2638.09095
$ws.Range("L109").Value = 10081.5792
$ws.Range("M109").Value = -1598.09095
$ws.Range("N109").Value = -12161.5792

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 940.9091
$ws.Range("I113").Value = 499.66666
$ws.Range("J113").Value = 1246.3846
$ws.Range("K113").Value = 1498.99998
$ws.Range("L113").Value = 3739.1538
$ws.Range("M113").Value = 671.0000199999999
$ws.Range("N113").Value = -8079.1538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 968.44446
$ws.Range("J121").Value = 1393.4
$ws.Range("L121").Value = 4180.200000000001
$ws.Range("N121").Value = -6800.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2016.175
$ws.Range("I102").Value = 2390.4644
$ws.Range("J102").Value = 1142.8334
$ws.Range("K102").Value = 2390.4644
$ws.Range("L102").Value = 1142.8334
$ws.Range("M102").Value = -768.4643999999998
$ws.Range("N102").Value = -4386.8334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3565.5667
$ws.Range("I126").Value = 2316.875
$ws.Range("J126").Value = 4992.643
$ws.Range("K126").Value = 6950.625
$ws.Range("L126").Value = 14977.929
$ws.Range("M126").Value = -4480.625
$ws.Range("N126").Value = -19917.929

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1458.5294
$ws.Range("I93").Value = 1416.6
$ws.Range("J93").Value = 1476
$ws.Range("K93").Value = 1416.6
$ws.Range("L93").Value = 1476
$ws.Range("M93").Value = -168.5999999999999
$ws.Range("N93").Value = -3972

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 38474104
$ws.Range("I136").Value = 62502548
$ws.Range("J136").Value = 28588.2
$ws.Range("K136").Value = 187507644
$ws.Range("L136").Value = 85764.60000000001
$ws.Range("M136").Value = -187505094
$ws.Range("N136").Value = -90864.60000000001
